$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Data table: each row is (ExcelRow, Col2..Col15 values for columns B..O)
# $null entries mean "leave unchanged" (columns E, K, N stay at 0)
$data = @(
    ,@(2, 1.079090931672567, 0.2813993398229684, 0.2247967378451534, $null, 1.138276000724971, 0.5663902846377198, 0.7088038347337928, 0.7023932395227011, 0.2119088042265709, $null, 0.3629779441536414, 0.2816805485701721, $null, 2.526337476135652)
    ,@(3, 0.9704600751207977, 0.2680210114067449, 0.2249199950571352, $null, 1.147957397497322, 0.571969731404316, 0.7155176558431435, 0.7132516209693822, 0.2136904709234564, $null, 0.3590213327849412, 0.2635630051467643, $null, 2.551955579266348)
    ,@(4, 0.9036174569251898, 0.2597901016218316, 0.2250656620387588, $null, 1.154568677939494, 0.5758328530807546, 0.7199814947120302, 0.7203469292444282, 0.2148510687524112, $null, 0.3567037669240278, 0.2524562886891744, $null, 2.569318038925871)
    ,@(5, 0.8763447266841808, 0.2564320699858058, 0.2251426778023671, $null, 1.157430545643642, 0.577516972446773, 0.7218864637856726, 0.7233459508868005, 0.2153408082382917, $null, 0.3557876057802218, 0.2479349894495897, $null, 2.576803833957555)
    ,@(6, 0.8718141298291471, 0.2558742464426871, 0.225156534211699, $null, 1.15791588714422, 0.5778032518251592, 0.7222079728632451, 0.7238504349257724, 0.2154231437769667, $null, 0.355637189531322, 0.2471845306564404, $null, 2.578071630677343)
    ,@(7, 0.9032497814995963, 0.2597448293060154, 0.2250666291389791, $null, 1.15460659496933, 0.5758551209954916, 0.7200068378592448, 0.720386939396283, 0.2148576055459355, $null, 0.3566912966069538, 0.2523952930287976, $null, 2.569417333258443)
    ,@(8, 1.04166587795379, 0.2767901022389196, 0.2248247448071297, $null, 1.141475755912047, 0.5682232513195444, 0.7110479003445747, 0.7060483020757022, 0.2125093068336827, $null, 0.3615905867436382, 0.2754301973912376, $null, 2.534831618540736)
    ,@(9, 1.311887316821924, 0.3100730921034938, 0.2249033801543661, $null, 1.121016584568345, 0.5567319123753691, 0.6961875735808505, 0.68133003331838, 0.2084318574471578, $null, 0.372079378404564, 0.3207266898869747, $null, 2.479971843317031)
    ,@(10, 1.509593734353189, 0.3344263052370025, 0.2252953225857652, $null, 1.109209725243055, 0.5504145006780092, 0.6869182958183586, 0.6652444124353618, 0.2057559971085734, $null, 0.3803159783233667, 0.3540669102491023, $null, 2.447579079107399)
    ,@(11, 1.599339348845433, 0.3454810096204994, 0.2255455815091096, $null, 1.104538546792611, 0.5480035395603124, 0.6830589527253608, 0.6583776663515764, 0.2046077549860534, $null, 0.3841768465659214, 0.3692444285024834, $null, 2.43456320591784)
    ,@(12, 1.633294230008119, 0.3496634695534908, 0.2256506420652542, $null, 1.102870313185242, 0.5471572470544004, 0.6816488644810477, 0.6558422726692825, 0.2041828432755484, $null, 0.3856551133287667, 0.3749929989887946, $null, 2.429881914440273)
    ,@(13, 1.6259828009903, 0.3487628721770761, 0.2256275584905936, $null, 1.103225121042769, 0.5473365433950619, 0.6819502681941287, 0.6563854275137935, 0.2042739155444107, $null, 0.3853360222960873, 0.3737548954902223, $null, 2.430879101284603)
    ,@(14, 1.602133447978588, 0.3458251793681768, 0.2255540189344813, $null, 1.104399283193786, 0.5479325773536914, 0.6829419146635161, 0.6581677769214771, 0.2045725989067932, $null, 0.3842981398679655, 0.3697173455684251, $null, 2.434173110885553)
    ,@(15, 1.587521079464864, 0.3440252645308135, 0.2255103126406084, $null, 1.105131597932051, 0.5483063536171713, 0.683556014977647, 0.6592679710706975, 0.2047568400503383, $null, 0.3836645172753634, 0.3672443709273878, $null, 2.436223030533711)
    ,@(16, 1.503724593054869, 0.3337033510073582, 0.2252804109244906, $null, 1.109529080549045, 0.5505813875249004, 0.6871777011959139, 0.6657022467312217, 0.2058324243767711, $null, 0.3800659423415311, 0.3530752092540439, $null, 2.448464323398937)
    ,@(17, 1.452267394816715, 0.3273649017425271, 0.2251577656645907, $null, 1.112406039997119, 0.5520956880706578, 0.6894909883345832, 0.669764958828809, 0.2065099222470028, $null, 0.3778874217305486, 0.3443853961586782, $null, 2.456414601959082)
    ,@(18, 1.422652632899485, 0.3237169814903496, 0.2250940000845958, $null, 1.114126661127706, 0.5530102328719195, 0.6908551615473328, 0.6721441387060789, 0.2069060989258151, $null, 0.3766451334097951, 0.3393882933036778, $null, 2.461149249650987)
    ,@(19, 1.412622569661437, 0.3224814879485507, 0.2250735757464639, $null, 1.11472054690848, 0.5533273594539097, 0.691322823963624, 0.6729569703721481, 0.2070413542354697, $null, 0.376226365050826, 0.3376965520623685, $null, 2.46278011191346)
    ,@(20, 1.457746974531403, 0.3280398716572961, 0.225170120496756, $null, 1.112092965385145, 0.5519299790538383, 0.6892412543496746, 0.6693280857360406, 0.2064371291048026, $null, 0.3781182186201875, 0.345310335996281, $null, 2.455551527827026)
    ,@(21, 1.609139408438864, 0.3466881542816793, 0.2255753403768637, $null, 1.10405167191152, 0.5477556968833852, 0.682649250271254, 0.6576424958046303, 0.2044845997843581, $null, 0.3846025515707794, 0.3709032431078469, $null, 2.433198860450091)
    ,@(22, 1.707908190966521, 0.3588540792550532, 0.2259001445559221, $null, 1.099382832635349, 0.5454163179422409, 0.6786403587270087, 0.6503835791043535, 0.2032662191191221, $null, 0.3889350250025245, 0.3876363484480905, $null, 2.420033044277659)
    ,@(23, 1.655210169259306, 0.3523629930125196, 0.225721321432097, $null, 1.10182100317904, 0.5466292757114815, 0.6807525936539633, 0.6542231574241306, 0.2039112186429985, $null, 0.3866140981123607, 0.3787050945590664, $null, 2.426927780980492)
    ,@(24, 1.455269753973084, 0.3277347297641029, 0.2251645138643354, $null, 1.112234298905832, 0.552004759199427, 0.6893540524476478, 0.6695254606539649, 0.2064700180932171, $null, 0.3780138436501659, 0.344892174446251, $null, 2.455941213041228)
    ,@(25, 1.238924001559099, 0.3010857731397039, 0.224823211958352, $null, 1.12598502581416, 0.5594679023665989, 0.6999180042384125, 0.6876527730892388, 0.2094786211230897, $null, 0.3691482357497762, 0.3084609926933908, $null, 2.493424042247185)
)

foreach ($row in $data) {
    $excelRow = $row[0]
    for ($i = 1; $i -lt $row.Count; $i++) {
        $val = $row[$i]
        if ($null -ne $val) {
            $col = $i + 1  # columns B..O map to column index 2..15
            $ws.Cells.Item($excelRow, $col).Value = $val
        }
    }
}